$p = $ppt.ActivePresentation
$s = $p.Slides.Item(37)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(5)
$r1 = $para.Runs(1)
$r2 = $para.Runs(2)
$r1.Text = "Das ganze gilt auch für die umgekehrten Fall mit den Nachfragern "
$r2.Text = ""
